$wb = $excel.ActiveWorkbook

# --- Metadata: bump "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 06:30 PM"

# --- Top Gainers: refresh rows 42-63 with the new snapshot ---
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Cells.Item(42, 2).Value = "INDOTHAI"
$gainers.Cells.Item(42, 3).Value = 4.8064
$gainers.Cells.Item(42, 4).Value = 4.5349
$gainers.Cells.Item(42, 5).Value = 43.748
$gainers.Cells.Item(43, 2).Value = "SANDUMA"
$gainers.Cells.Item(43, 3).Value = 4.593
$gainers.Cells.Item(43, 4).Value = 2.1405
$gainers.Cells.Item(43, 5).Value = 30.2813
$gainers.Cells.Item(44, 2).Value = "LLOYDSENT"
$gainers.Cells.Item(44, 3).Value = 4.5646
$gainers.Cells.Item(44, 4).Value = 1.8339
$gainers.Cells.Item(44, 5).Value = 11.234
$gainers.Cells.Item(45, 2).Value = "STAR"
$gainers.Cells.Item(45, 3).Value = 4.5025
$gainers.Cells.Item(45, 4).Value = 4.4319
$gainers.Cells.Item(45, 5).Value = 3.662
$gainers.Cells.Item(46, 2).Value = "RECLTD"
$gainers.Cells.Item(46, 3).Value = 4.4992
$gainers.Cells.Item(46, 4).Value = 3.4756
$gainers.Cells.Item(46, 5).Value = 3.4062
$gainers.Cells.Item(47, 2).Value = "NBCC"
$gainers.Cells.Item(47, 3).Value = 4.4511
$gainers.Cells.Item(47, 4).Value = 3.1605
$gainers.Cells.Item(47, 5).Value = 7.6018
$gainers.Cells.Item(48, 2).Value = "GPPL"
$gainers.Cells.Item(48, 3).Value = 4.4154
$gainers.Cells.Item(48, 4).Value = 3.4073
$gainers.Cells.Item(48, 5).Value = 5.0497
$gainers.Cells.Item(49, 2).Value = "BIL"
$gainers.Cells.Item(49, 3).Value = 4.3654
$gainers.Cells.Item(49, 4).Value = 9.122199999999999
$gainers.Cells.Item(49, 5).Value = -0.3203
$gainers.Cells.Item(50, 2).Value = "HUDCO"
$gainers.Cells.Item(50, 3).Value = 4.3201
$gainers.Cells.Item(50, 4).Value = 3.8924
$gainers.Cells.Item(50, 5).Value = 5.3884
$gainers.Cells.Item(51, 2).Value = "SGMART"
$gainers.Cells.Item(51, 3).Value = 4.2736
$gainers.Cells.Item(51, 4).Value = 8.258900000000001
$gainers.Cells.Item(51, 5).Value = 2.5381
$gainers.Cells.Item(52, 2).Value = "MRPL"
$gainers.Cells.Item(52, 3).Value = 4.2642
$gainers.Cells.Item(52, 4).Value = 9.7103
$gainers.Cells.Item(52, 5).Value = 20.0542
$gainers.Cells.Item(53, 2).Value = "JKIL"
$gainers.Cells.Item(53, 3).Value = 4.1372
$gainers.Cells.Item(53, 4).Value = 2.9463
$gainers.Cells.Item(53, 5).Value = 1.7584
$gainers.Cells.Item(54, 2).Value = "SAMBHV"
$gainers.Cells.Item(54, 3).Value = 4.1349
$gainers.Cells.Item(54, 4).Value = 2.624
$gainers.Cells.Item(54, 5).Value = 5.167
$gainers.Cells.Item(55, 2).Value = "SAPPHIRE"
$gainers.Cells.Item(55, 3).Value = 4.1265
$gainers.Cells.Item(55, 4).Value = 1.7633
$gainers.Cells.Item(55, 5).Value = -0.7999000000000001
$gainers.Cells.Item(56, 2).Value = "PVRINOX"
$gainers.Cells.Item(56, 3).Value = 4.1118
$gainers.Cells.Item(56, 4).Value = 6.2102
$gainers.Cells.Item(56, 5).Value = 14.707
$gainers.Cells.Item(57, 2).Value = "KERNEX"
$gainers.Cells.Item(57, 3).Value = 4.0782
$gainers.Cells.Item(57, 4).Value = 7.542
$gainers.Cells.Item(57, 5).Value = 27.2033
$gainers.Cells.Item(58, 2).Value = "SUNFLAG"
$gainers.Cells.Item(58, 3).Value = 3.997
$gainers.Cells.Item(58, 4).Value = 4.333
$gainers.Cells.Item(58, 5).Value = 4.6312
$gainers.Cells.Item(59, 2).Value = "CMSINFO"
$gainers.Cells.Item(59, 3).Value = 3.9096
$gainers.Cells.Item(59, 4).Value = 2.6872
$gainers.Cells.Item(59, 5).Value = 2.8935
$gainers.Cells.Item(60, 2).Value = "GMBREW"
$gainers.Cells.Item(60, 3).Value = 3.8999
$gainers.Cells.Item(60, 4).Value = -0.53
$gainers.Cells.Item(60, 5).Value = 79.029
$gainers.Cells.Item(61, 2).Value = "GREENLAM"
$gainers.Cells.Item(61, 3).Value = 3.8946
$gainers.Cells.Item(61, 4).Value = 3.5858
$gainers.Cells.Item(61, 5).Value = 10.721
$gainers.Cells.Item(62, 2).Value = "APARINDS"
$gainers.Cells.Item(62, 3).Value = 3.8924
$gainers.Cells.Item(62, 4).Value = 8.3414
$gainers.Cells.Item(62, 5).Value = 15.5876
$gainers.Cells.Item(63, 2).Value = "HITECHGEAR"
$gainers.Cells.Item(63, 3).Value = 3.8587
$gainers.Cells.Item(63, 4).Value = 1.1486
$gainers.Cells.Item(63, 5).Value = 9.9254

# --- Top Losers: refresh rows 71-76 with the new snapshot ---
$losers = $wb.Worksheets.Item("Top Losers")
$losers.Cells.Item(71, 2).Value = "FCL"
$losers.Cells.Item(71, 3).Value = -2.3453
$losers.Cells.Item(71, 4).Value = -2.616
$losers.Cells.Item(71, 5).Value = -0.02
$losers.Cells.Item(72, 2).Value = "DEEDEV"
$losers.Cells.Item(72, 3).Value = -2.3136
$losers.Cells.Item(72, 4).Value = -6.6339
$losers.Cells.Item(72, 5).Value = -7.4039
$losers.Cells.Item(73, 2).Value = "WEALTH"
$losers.Cells.Item(73, 3).Value = -2.3047
$losers.Cells.Item(73, 4).Value = -3.8606
$losers.Cells.Item(73, 5).Value = -2.8234
$losers.Cells.Item(74, 2).Value = "RATNAMANI"
$losers.Cells.Item(74, 3).Value = -2.2788
$losers.Cells.Item(74, 4).Value = -0.4626
$losers.Cells.Item(74, 5).Value = 0.8712
$losers.Cells.Item(75, 2).Value = "CSBBANK"
$losers.Cells.Item(75, 3).Value = -2.2695
$losers.Cells.Item(75, 4).Value = 2.3137
$losers.Cells.Item(75, 5).Value = 10.6999
$losers.Cells.Item(76, 2).Value = "BBOX"
$losers.Cells.Item(76, 3).Value = -2.2639
$losers.Cells.Item(76, 4).Value = -4.7636
$losers.Cells.Item(76, 5).Value = 5.1528

# --- distance from Dma50: drop the trailing N/A rows (31-33) ---
$dma = $wb.Worksheets.Item("distance from Dma50")
$dma.Rows("31:33").Delete()

Write-Output "edit complete"